$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 50. This shifts the existing rows 50-76
# down to 51-77, preserving all of their data untouched.
$ws.Rows(50).Insert()

# Populate the newly inserted row 50 with the latest weekly price record.
$ws.Cells.Item(50, 1).Value = 4
$ws.Cells.Item(50, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(50, 3).Value = "Los Lagos"
$ws.Cells.Item(50, 4).Value = 44523
$ws.Cells.Item(50, 5).Value = 10
$ws.Cells.Item(50, 6).Value = 100112052
$ws.Cells.Item(50, 7).Value = "Albahaca"
$ws.Cells.Item(50, 8).Value = "Sin especificar"
$ws.Cells.Item(50, 9).Value = "Primera"
$ws.Cells.Item(50, 10).Value = 90
$ws.Cells.Item(50, 11).Value = 8000
$ws.Cells.Item(50, 12).Value = 8000
$ws.Cells.Item(50, 13).Value = 8000
$ws.Cells.Item(50, 14).Value = "$/docena de matas"
$ws.Cells.Item(50, 15).Value = "Región Metropolitana"
$ws.Cells.Item(50, 16).Value = 1333
$ws.Cells.Item(50, 17).Value = 6
$ws.Cells.Item(50, 18).Value = "Hortaliza"
